# 24 Feb 2025 Part 1
# Add a new team member "zuhran.ahmed@ascend.com.sa" to the Projects sheet.
# The new record is inserted as a new row directly above the existing
# "athar.ali@ascend.com.sa" row (row 43), pushing that row and all rows
# below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 43 (shifts rows 43:47 down to 44:48).
$ws.Rows("43:43").Insert()

$ws.Cells.Item(43, 1).Value = "zuhran.ahmed@ascend.com.sa"
$ws.Cells.Item(43, 2).Value = "Approved"
$ws.Cells.Item(43, 3).Value = "MoH Compliance Program_Digital  /  Digital Innovation - EFX - Inspection`nN/A - Digital Innovation-Investment Work-Non-PO  /  Astrom (NEC-MOHU-Others)`nN/A - Digital Innovation-Investment Work-Non-PO  /  No Tasks"
